# Apply the edits described by the commit diff to Sheet1 of the workbook:
#   - tweak several input numbers in the second ("理想") table (rows 10-15)
#   - append a brand-new category row (16) called "信息" with B/C inputs
#     and the two delta formulas that Excel auto-filled for it
#   - move the active selection to C14 to match the saved view state
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# --- Update existing input values in the second table (rows 10-15) ---
$ws.Range("C10").Value = 1.2
$ws.Range("D10").Value = 1.3

$ws.Range("D11").Value = 1.3

$ws.Range("D12").Value = 1.5

$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 1.6

$ws.Range("C14").Value = 2.5
$ws.Range("D14").Value = 1.7

$ws.Range("C15").Value = 2.5

# --- Add a new row (16) for a new category "信息" ---
# Copy the header-cell formatting of A1 (the bold/compat-font style) onto A16
# so it matches the style used for the other "s=1" category label.
$ws.Range("A1").Copy()
$ws.Range("A16").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A16").Value = "信息"

$ws.Range("B16").Value = 2.5
$ws.Range("C16").Value = 1.8

$ws.Range("J16").Formula = "=B16-B15"
$ws.Range("K16").Formula = "=C16-C15"

# --- Update the active selection to match the final saved state ---
$null = $ws.Range("C14").Select()
